$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-10-05 to 2023-10-08
# (date serial numbers 45204 -> 45207), preserving existing cell formatting.
$ws.Range("C2:C5").Value = 45207
